$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 81.64706
$ws.Range("I5").Value = 52.57143
$ws.Range("K5").Value = 52.57143
$ws.Range("M5").Value = 62.42857

$ws.Range("H9").Value = 217.625
$ws.Range("I9").Value = 253.36363
$ws.Range("J9").Value = 139
$ws.Range("K9").Value = 253.36363
$ws.Range("L9").Value = 139
$ws.Range("M9").Value = -84.36363
$ws.Range("N9").Value = -477

$ws.Range("H39").Value = 434.25
$ws.Range("I39").Value = 379
$ws.Range("J39").Value = 600
$ws.Range("K39").Value = 1137
$ws.Range("L39").Value = 1800
$ws.Range("M39").Value = -841
$ws.Range("N39").Value = -2392

$ws.Range("H40").Value = 4340.4517
$ws.Range("I40").Value = 3106.28
$ws.Range("K40").Value = 3106.28
$ws.Range("M40").Value = -2931.28

$ws.Range("H135").Value = 496.06668
$ws.Range("I135").Value = 496.06668
$ws.Range("K135").Value = 4464.60012
$ws.Range("M135").Value = -1929.60012

$ws.Range("H137").Value = 2322.6924
$ws.Range("I137").Value = 1044.037
$ws.Range("J137").Value = 3703.64
$ws.Range("K137").Value = 3132.111
$ws.Range("L137").Value = 11110.92
$ws.Range("M137").Value = -582.1109999999999
$ws.Range("N137").Value = -16210.92

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 208.125
$ws.Range("I4").Value = 164.33333
$ws.Range("J4").Value = 339.5
$ws.Range("K4").Value = 164.33333
$ws.Range("L4").Value = 339.5
$ws.Range("M4").Value = -48.33332999999999
$ws.Range("N4").Value = -571.5

$ws.Range("H32").Value = 2145.3333
$ws.Range("I32").Value = 936.3125
$ws.Range("K32").Value = 936.3125
$ws.Range("M32").Value = -649.3125

$ws.Range("H110").Value = 860.9286
$ws.Range("I110").Value = 858.0833
$ws.Range("J110").Value = 878
$ws.Range("K110").Value = 858.0833
$ws.Range("L110").Value = 878
$ws.Range("M110").Value = 1186.9167
$ws.Range("N110").Value = -4968

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5729
$ws.Range("I107").Value = 4888.6665
$ws.Range("K107").Value = 4888.6665
$ws.Range("M107").Value = -2968.6665

$ws.Range("H134").Value = 2499.5
$ws.Range("I134").Value = 2499.5
$ws.Range("K134").Value = 7498.5
$ws.Range("M134").Value = -4963.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4750
$ws.Range("I25").Value = 4500
$ws.Range("K25").Value = 4500
$ws.Range("M25").Value = -4326

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 722.0345
$ws.Range("I5").Value = 717.4545000000001
$ws.Range("J5").Value = 724.8333
$ws.Range("K5").Value = 2152.3635
$ws.Range("L5").Value = 2174.4999
$ws.Range("M5").Value = -2040.3635
$ws.Range("N5").Value = -2398.4999

$ws.Range("H37").Value = 163873
$ws.Range("J37").Value = 163873
$ws.Range("L37").Value = 491619
$ws.Range("N37").Value = -491843

$ws.Range("H68").Value = 2799.4
$ws.Range("I68").Value = 2748.5
$ws.Range("K68").Value = 8245.5
$ws.Range("M68").Value = -7434.5

$ws.Range("H71").Value = 2799.4
$ws.Range("I71").Value = 2748.5
$ws.Range("K71").Value = 24736.5
$ws.Range("M71").Value = -20680.5

$ws.Range("H107").Value = 351.5
$ws.Range("I107").Value = 351.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1054.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 865.5
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 870.125
$ws.Range("I113").Value = 662.4
$ws.Range("J113").Value = 1216.3334
$ws.Range("K113").Value = 1987.2
$ws.Range("L113").Value = 3649.0002
$ws.Range("M113").Value = 182.8000000000002
$ws.Range("N113").Value = -7989.0002

$ws.Range("H121").Value = 243
$ws.Range("I121").Value = 157.33333
$ws.Range("J121").Value = 500
$ws.Range("K121").Value = 471.99999
$ws.Range("L121").Value = 1500
$ws.Range("M121").Value = 838.00001
$ws.Range("N121").Value = -4120

$ws.Range("H135").Value = 722.0345
$ws.Range("I135").Value = 717.4545000000001
$ws.Range("J135").Value = 724.8333
$ws.Range("K135").Value = 6457.0905
$ws.Range("L135").Value = 6523.4997
$ws.Range("M135").Value = -3922.0905
$ws.Range("N135").Value = -11593.4997

$ws.Range("H140").Value = 4216
$ws.Range("I140").Value = 4020
$ws.Range("K140").Value = 12060
$ws.Range("M140").Value = -6880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 10000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -10504

$ws.Range("H132").Value = 5228.353
$ws.Range("I132").Value = 3838.2
$ws.Range("J132").Value = 7214.2856
$ws.Range("K132").Value = 11514.6
$ws.Range("L132").Value = 21642.8568
$ws.Range("M132").Value = -8984.599999999999
$ws.Range("N132").Value = -26702.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4370
$ws.Range("I22").Value = 3999.5
$ws.Range("K22").Value = 3999.5
$ws.Range("M22").Value = -3704.5

$ws.Range("H27").Value = 4370
$ws.Range("I27").Value = 3999.5
$ws.Range("K27").Value = 3999.5
$ws.Range("M27").Value = -3892.5

$ws.Range("H32").Value = 2180
$ws.Range("I32").Value = 1766.6666
$ws.Range("K32").Value = 1766.6666
$ws.Range("M32").Value = -1449.6666

$ws.Range("H46").Value = 5792.778
$ws.Range("I46").Value = 4330.3
$ws.Range("J46").Value = 7620.875
$ws.Range("K46").Value = 4330.3
$ws.Range("L46").Value = 7620.875
$ws.Range("M46").Value = -4142.3
$ws.Range("N46").Value = -7996.875

$ws.Range("H61").Value = 2371.6843
$ws.Range("I61").Value = 1337.5333
$ws.Range("K61").Value = 1337.5333
$ws.Range("M61").Value = -1135.5333

$ws.Range("H68").Value = 6599.7
$ws.Range("I68").Value = 2999
$ws.Range("J68").Value = 8142.857
$ws.Range("K68").Value = 2999
$ws.Range("L68").Value = 8142.857
$ws.Range("M68").Value = -2250
$ws.Range("N68").Value = -9640.857

$ws.Range("H71").Value = 6599.7
$ws.Range("I71").Value = 2999
$ws.Range("J71").Value = 8142.857
$ws.Range("K71").Value = 14995
$ws.Range("L71").Value = 40714.285
$ws.Range("M71").Value = -11251
$ws.Range("N71").Value = -48202.285

$ws.Range("H93").Value = 945
$ws.Range("I93").Value = 890
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 890
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 358
$ws.Range("N93").Value = -3496

$ws.Range("H113").Value = 2371.6843
$ws.Range("I113").Value = 1337.5333
$ws.Range("K113").Value = 1337.5333
$ws.Range("M113").Value = 832.4666999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11916.667
$ws.Range("I62").Value = 11500
$ws.Range("K62").Value = 11500
$ws.Range("M62").Value = -10876

$ws.Range("H65").Value = 11916.667
$ws.Range("I65").Value = 11500
$ws.Range("K65").Value = 57500
$ws.Range("M65").Value = -54380

$ws.Range("H81").Value = 6396.727
$ws.Range("I81").Value = 2733.625
$ws.Range("J81").Value = 16165
$ws.Range("K81").Value = 5467.25
$ws.Range("L81").Value = 32330
$ws.Range("M81").Value = -4406.25
$ws.Range("N81").Value = -34452

$ws.Range("H84").Value = 6396.727
$ws.Range("I84").Value = 2733.625
$ws.Range("J84").Value = 16165
$ws.Range("K84").Value = 27336.25
$ws.Range("L84").Value = 161650
$ws.Range("M84").Value = -22032.25
$ws.Range("N84").Value = -172258
